# Update header text in D1: "E" -> "Error"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Error"

# Update the secant-method data table (rows 2-10, columns A-D)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1.00001

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = -0.281718171540955
$ws.Range("D3").Value = 1.00001

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 0.780202717105698
$ws.Range("C4").Value = -0.158693619249085
$ws.Range("D4").Value = 0.219797282894302

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 0.496678610138336
$ws.Range("C5").Value = 0.153218478153997
$ws.Range("D5").Value = 0.283524106967362

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 0.635952247009725
$ws.Range("C6").Value = -0.0190368325556636
$ws.Range("D6").Value = 0.139273636871389

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 0.62056039120211
$ws.Range("C7").Value = -0.0017111128981504
$ws.Range("D7").Value = 0.0153918558076145

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 0.61904026946867
$ws.Range("C8").Value = 0.0000240192824358942
$ws.Range("D8").Value = 0.0015201217334402

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 0.619061312380364
$ws.Range("C9").Value = -0.0000000293068544987563
$ws.Range("D9").Value = 0.0000210429116933586

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 0.6190612867363831
$ws.Range("C10").Value = -0.000000000000500710584105946
$ws.Range("D10").Value = 0.0000000256439803836273
